$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for D/E columns so numeric-looking strings (e.g. "1.002")
# are kept as text instead of being auto-converted to numbers by Excel.
$numRange = $ws.Range("D2:E51")
$numRange.NumberFormat = "@"

$ws.Range('D2').Value = '27.232.70'
$ws.Range('E2').Value = '  -2.06%  '
$ws.Range('D3').Value = '1.873.57'
$ws.Range('E3').Value = '  -1.51%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '307.89'
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('D7').Value = '0.5181'
$ws.Range('E7').Value = '  +3.28%  '
$ws.Range('D8').Value = '0.3761'
$ws.Range('E8').Value = '  -1.25%  '
$ws.Range('D9').Value = '0.07155'
$ws.Range('E9').Value = '  -1.67%  '
$ws.Range('D10').Value = '20.89'
$ws.Range('E10').Value = '  +0.45%  '
$ws.Range('D11').Value = '0.8860'
$ws.Range('E11').Value = '  -2.40%  '
$ws.Range('D12').Value = '1.880.49'
$ws.Range('E12').Value = '  -1.38%  '
$ws.Range('D13').Value = '0.07595'
$ws.Range('E13').Value = '  -0.74%  '
$ws.Range('E14').Value = '  -2.43%  '
$ws.Range('D15').Value = '89.53'
$ws.Range('E16').Value = '  -0.18%  '
$ws.Range('D17').Value = '0.000008572'
$ws.Range('E17').Value = '  -1.38%  '
$ws.Range('E18').Value = '  -2.18%  '
$ws.Range('E19').Value = '  -0.18%  '
$ws.Range('D20').Value = '27.280.32'
$ws.Range('E20').Value = '  -2.03%  '
$ws.Range('D21').Value = '5.048'
$ws.Range('E21').Value = '  -2.19%  '
$ws.Range('D22').Value = '2.133.00'
$ws.Range('E22').Value = '  +1.53%  '
$ws.Range('E23').Value = '  -1.50%  '
$ws.Range('D24').Value = '6.482'
$ws.Range('E24').Value = '  -1.69%  '
$ws.Range('D25').Value = '151.81'
$ws.Range('E25').Value = '  -1.56%  '
$ws.Range('D26').Value = '1.853'
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = '2.187'
$ws.Range('E27').Value = '  -1.80%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '18.07'
$ws.Range('E28').Value = '  -1.62%  '
$ws.Range('D29').Value = '113.11'
$ws.Range('E29').Value = '  -1.83%  '
$ws.Range('E30').Value = '  -2.99%  '
$ws.Range('D31').Value = '4.714'
$ws.Range('E31').Value = '  +1.71%  '
$ws.Range('D32').Value = '0.09051'
$ws.Range('E32').Value = '  +0.97%  '
$ws.Range('D33').Value = '0.05187'
$ws.Range('E33').Value = '  -1.23%  '
$ws.Range('D34').Value = '3.088'
$ws.Range('E34').Value = '  -3.74%  '
$ws.Range('D35').Value = '0.7608'
$ws.Range('E35').Value = '  -0.64%  '
$ws.Range('D36').Value = '1.182'
$ws.Range('E36').Value = '  -3.93%  '
$ws.Range('D37').Value = '0.02055'
$ws.Range('E37').Value = '  -0.12%  '
$ws.Range('D38').Value = '2.556'
$ws.Range('E38').Value = '  +0.45%  '
$ws.Range('D39').Value = '3.045'
$ws.Range('E39').Value = '  +0.95%  '
$ws.Range('D40').Value = '1.084'
$ws.Range('E40').Value = '  -1.11%  '
$ws.Range('D41').Value = '0.5469'
$ws.Range('E41').Value = '  -1.38%  '
$ws.Range('D42').Value = '6.695'
$ws.Range('E42').Value = '  -3.97%  '
$ws.Range('D43').Value = '115.69'
$ws.Range('E43').Value = '  +4.18%  '
$ws.Range('D44').Value = '8.583'
$ws.Range('E44').Value = '  +1.22%  '
$ws.Range('E45').Value = '  -2.08%  '
$ws.Range('D46').Value = '0.4713'
$ws.Range('E46').Value = '  -1.49%  '
$ws.Range('D47').Value = '10.21'
$ws.Range('E47').Value = '  -3.75%  '
$ws.Range('D48').Value = '1.002'
$ws.Range('E48').Value = '  -0.16%  '
$ws.Range('D49').Value = '1.582'
$ws.Range('E49').Value = '  -2.92%  '
$ws.Range('D50').Value = '65.31'
$ws.Range('E50').Value = '  -2.89%  '
$ws.Range('D51').Value = '36.60'
$ws.Range('E51').Value = '  -1.18%  '

# Restore default (unstyled) formatting for the D/E columns so that cell
# styling matches the original workbook (no explicit style index).
$numRange.Style = "Normal"
